# Update the "Productos" sheet: a new product ("Tela Azul", id 3) is
# inserted as row 4, and the old row-4 product ("Tela Azul", id 6) slides
# down to row 5, replacing the previous "Tela Verde" row which is removed.
#
# Before:
#   Row4: 6 | Tela Azul  | 8  | 8000.0  | 64000.0
#   Row5: 7 | Tela Verde | 11 | 15000.0 | 165000.0
#
# After:
#   Row4: 3 | Tela Azul  | 3  | 5000.0  | 15000.0
#   Row5: 6 | Tela Azul  | 8  | 8000.0  | 64000.0

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores every value (including numeric-looking ones like "6" or
# "8000.0") as text rather than real numbers. Mark the numeric-looking cells
# we are about to rewrite as Text first so Excel keeps the new values as
# text too, instead of silently converting them into real numbers.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("C4:E4").NumberFormat = "@"
$ws.Range("A5").NumberFormat = "@"
$ws.Range("C5:E5").NumberFormat = "@"

# Row 4 -> new product "Tela Azul" (id 3); B4 already reads "Tela Azul".
$ws.Cells.Item(4, 1).Value = "3"
$ws.Cells.Item(4, 3).Value = "3"
$ws.Cells.Item(4, 4).Value = "5000.0"
$ws.Cells.Item(4, 5).Value = "15000.0"

# Row 5 -> previous "Tela Azul" (id 6) row, "Tela Verde" row is dropped
$ws.Cells.Item(5, 1).Value = "6"
$ws.Cells.Item(5, 2).Value = "Tela Azul"
$ws.Cells.Item(5, 3).Value = "8"
$ws.Cells.Item(5, 4).Value = "8000.0"
$ws.Cells.Item(5, 5).Value = "64000.0"
